# Commit: Fri, Jul 10, 2020 8:05:57 AM
#
# The only semantically reachable change via the PowerPoint object model is
# the table style applied to the single table on slide 16: it moves from the
# deck's custom "Table_0" style ({13D1B01A-9CEB-476D-B096-123662127712}) to
# the built-in table style {A760E8D9-90E8-443B-A639-3FE483DBB335}.
#
# (Table.Style is read-only as a plain property in this object model --
# attempting `$tbl.Style = "{...}"` raises "Table styles cannot be assigned
# through a property - call Table.ApplyStyle("{GUID}") instead", which is
# exactly what we do below.)

$p = $ppt.ActivePresentation

$targetStyleId = "{A760E8D9-90E8-443B-A639-3FE483DBB335}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}
